$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Okay." -> "Thank you." (two occurrences: after C., after F.)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Okay.", $true, $false, $false, $false, $false, $true, 1, $false, "Thank you.", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Replace the first "Response needed." paragraph (currently right
#    after the "...C only is sequestered in trees, not CO2." remark)
#    with the new bold reply about sequestration wording. The
#    "D. Appropriate use of statistics..." remark and the
#    "Response needed." paragraph that follows it are left untouched.
# ------------------------------------------------------------------

$p = $d.Paragraphs(14)
$r = $p.Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = ""

$ip = $d.Paragraphs(14).Range
$ip.MoveEnd(1, -1) | Out-Null
$ip.Collapse(1)

$ip.InsertAfter('We have replaced “CO')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Collapse(0)

$ip.InsertAfter('2')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Font.Subscript = $true
$ip.Collapse(0)

$ip.InsertAfter(' sequestration” or “C sequestration” with “CO')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Collapse(0)

$ip.InsertAfter('2')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Font.Subscript = $true
$ip.Collapse(0)

$ip.InsertAfter(' uptake” when referring to ecosystem-atmosphere exchange of CO')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Collapse(0)

$ip.InsertAfter('2')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Font.Subscript = $true
$ip.Collapse(0)

$ip.InsertAfter(' and use “sequestration” only to refer to woody growth. Further, we have removed all use of “CO2 sequestration” and now use only “C sequestration”.')
$ip.Font.Bold = $true
$ip.Font.BoldBi = $true
$ip.Collapse(0)

Write-Output "done"
